# "Moved structural types down a slot."
#
# The block of "structural" one-hot rows (list / map / reserved / list-of-lists /
# list-of-maps / reserved) that used to sit in rows 20-26 using columns F/G/H as the
# one-hot selector now sits one "slot" earlier, in rows 12-18, using columns E/F/G as
# the one-hot selector. A brand-new "list-of" entry is inserted into that block (row 15).
# Rows 20-26 become plain "reserved" rows with every bit cleared to empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 12-18: new structural-type block (E/F/G one-hot), column A label updated ---

# Row 12: list -> 1,0,0 in E/F/G
$ws.Range("A12").Value = "list"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0

# Row 13: map -> 0,1,0 in E/F/G
$ws.Range("A13").Value = "map"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0

# Row 14: reserved -> 1,1,0 in E/F/G
$ws.Range("A14").Value = "reserved"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0

# Row 15: list-of (new entry) -> 0,0,1 in E/F/G
$ws.Range("A15").Value = "list-of"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1

# Row 16: list-of-lists -> 1,0,1 in E/F/G
$ws.Range("A16").Value = "list-of-lists"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 1

# Row 17: list-of-maps -> 0,1,1 in E/F/G
$ws.Range("A17").Value = "list-of-maps"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1

# Row 18: reserved -> 1,1,1 in E/F/G
$ws.Range("A18").Value = "reserved"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1

# Row 19: reserved, one-hot moves to H (bit 7) instead of B/C/D/E being all set
$ws.Range("A19").Value = "reserved"
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("H19").Value = 1

# --- Rows 20-26: old structural-type block, now cleared down to plain "reserved" rows ---

$ws.Range("A20").Value = "reserved"
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""

$ws.Range("A21").Value = "reserved"
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = ""

$ws.Range("A22").Value = "reserved"
$ws.Range("F22").Value = ""
$ws.Range("G22").Value = ""
$ws.Range("H22").Value = ""

$ws.Range("A23").Value = "reserved"
$ws.Range("F23").Value = ""
$ws.Range("G23").Value = ""
$ws.Range("H23").Value = ""

$ws.Range("A24").Value = "reserved"
$ws.Range("F24").Value = ""
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = ""

$ws.Range("A25").Value = "reserved"
$ws.Range("F25").Value = ""
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = ""

$ws.Range("A26").Value = "reserved"
$ws.Range("F26").Value = ""
$ws.Range("G26").Value = ""
$ws.Range("H26").Value = ""

# --- Row 6: O6's lookup formula now references the relocated "list" row (J12) ---
$ws.Range("O6").Formula = "=J7+J12"

# --- Selection moved from G26 to O7 ---
$ws.Range("O7").Select()
